# Basic Bringup.docx edit script
# Reproduces the target diff: proofing-pass run splits (w:proofErr markers),
# a couple of xml:space="preserve" additions, and relocation of the
# "_GoBack" bookmark from the "Suggested Silver brazing" bullet to the
# title paragraph (with a new leading space run).

$d = $word.ActiveDocument

function Get-ParagraphsByText($needle) {
    $result = @()
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            $result += $p
        }
    }
    return $result
}

function Replace-ParagraphXml($needle, $xml) {
    $paras = Get-ParagraphsByText($needle)
    foreach ($p in $paras) {
        $p.Range.InsertXML($xml)
    }
}

# ---------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark off the "Suggested Silver brazing" bullet
#    (it will be re-created at the title once that edit is made below).
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. Title paragraph: add a leading space run, split "Bringup"/"Cryorefrigerator"
#    runs with spellcheck proofErr markers, and drop the new _GoBack bookmark
#    right after the new leading space.
# ---------------------------------------------------------------------
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Basic </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bringup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cryorefrigerator</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(1).Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 3. Introduction paragraph - gramStart/gramEnd + spellStart/spellEnd splits.
# ---------------------------------------------------------------------
$introXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">This document is an ongoing document detailing the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>current status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> of the quick test </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cryorefrigerator</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. The intent of this document is to note current problems and attempted solutions as well as providing a history for the project. </w:t></w:r></w:p>'
Replace-ParagraphXml "This document is an ongoing document" $introXml

# ---------------------------------------------------------------------
# 4. Goals paragraph - gramStart/gramEnd split around "quickly".
# ---------------------------------------------------------------------
$goalsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The goal of this project is to make a cold chamber that can </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>quickly</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and cost effectively reach</w:t></w:r><w:r><w:t xml:space="preserve"> ~1K </w:t></w:r><w:r><w:t>temperatures. The use for this system would be quick tests of samples prior to placing them in the main cryostats that take ages to reach temperature.</w:t></w:r></w:p>'
Replace-ParagraphXml "The goal of this project" $goalsXml

# ---------------------------------------------------------------------
# 5. "PT415 spira and o-ring groove..." bullet - spellStart/spellEnd splits.
# ---------------------------------------------------------------------
$pt415Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t xml:space="preserve">PT415 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t>spira</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t>o-ring</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> groove. Use recommended sizes based on the plate bottoming out and the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t>spira</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> being on the outside of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t>oring</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Replace-ParagraphXml "PT415 spira" $pt415Xml

# ---------------------------------------------------------------------
# 6. "55k stage will be made out of aluminum" bullet - gramStart/gramEnd split.
# ---------------------------------------------------------------------
$stage55kXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t xml:space="preserve">55k stage will be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t>made out of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> aluminum</w:t></w:r></w:p>'
Replace-ParagraphXml "55k stage will be made" $stage55kXml

# ---------------------------------------------------------------------
# 7. "Inventor professional..." paragraph - gramStart/gramEnd splits around
#    "tool" and "stainless steel".
# ---------------------------------------------------------------------
$inventorXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Inventor professional has a built in FEM solver for stress. Using this </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>tool</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> it is straightforward to get a rough idea of the deflection which the plate will see under vacuum. A pressure of 15psi was applied to the surface of the plate with the vacuum jacket marked as an immovable object. The result of the calculations is shown in the image below. The takeaway is that for a &#190;&#8221; </w:t></w:r><w:r><w:t xml:space="preserve">austenitic </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>stainless steel</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> plate a maximal deflection of 0.002&#8221; is expected.</w:t></w:r></w:p>'
Replace-ParagraphXml "Inventor professional" $inventorXml

# ---------------------------------------------------------------------
# 8. "A further test with gravity..." paragraph - gramStart/gramEnd split
#    around "have  a".
# ---------------------------------------------------------------------
$gravityXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">A further test with gravity found no difference. Also initial investigations into the natural harmonics of the steel plate show that modifying the orientation of the holes does not </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>have  a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> significant effect on the dominant vibrational mode. You can see that the spherical harmonics for a disk are predominant (as expected) with additional modes due to the corners vibrating. </w:t></w:r></w:p>'
Replace-ParagraphXml "A further test with gravity" $gravityXml

# ---------------------------------------------------------------------
# 9. "Spria-shield is a product..." paragraph - spellStart/spellEnd split
#    around "Spria".
# ---------------------------------------------------------------------
$spriaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Spria</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-shield is a product manufactured by Spira manufacturing. This product looks like a long metal spiral made of a flat continuous strip of metal. The intended use of Spira is providing EMI shielding by sandwiching the Spira between two metal surfaces. The Spira then provides conductive contact between the two metal surfaces screening EMI.</w:t></w:r></w:p>'
Replace-ParagraphXml "Spria-shield is a product" $spriaXml

# ---------------------------------------------------------------------
# 10. "The purpose of this is to reduce..." paragraph - spellStart/spellEnd
#     split around "fulled".
# ---------------------------------------------------------------------
$purposeXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The purpose of this is to reduce the radiated heat leak into the coldest part of the cryostat. (Pg.10 &#8220;Experimental techniques in condensed matter&#8221;). The first design question addressed was whether to use some pre-existing shields or to manufacture our own. We have decided to repurpose some shields which were laying around, whether this is the best option has not been </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fulled</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> explored. Of primary concern is the weight of the shields, with secondary concern that they consume too much of the space in the cryostat.</w:t></w:r></w:p>'
Replace-ParagraphXml "The purpose of this is to reduce" $purposeXml

# ---------------------------------------------------------------------
# 11. Both "Weight (lbs)" table header cells - spellStart/spellEnd split
#     around "lbs".
# ---------------------------------------------------------------------
$weightLbsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Weight (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lbs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>'
Replace-ParagraphXml "Weight (lbs)" $weightLbsXml

# ---------------------------------------------------------------------
# 12. "...As a result they will have to be cut to length..." paragraph -
#     gramStart/gramEnd split around "result".
# ---------------------------------------------------------------------
$cutLengthXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Also note that the radiation shields available are not the right </w:t></w:r><w:r><w:t>height</w:t></w:r><w:r><w:t xml:space="preserve"> for this purpose.</w:t></w:r><w:r><w:t xml:space="preserve"> As a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>result</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> they will have to be cut to length, which might be as much work as making custom shields.</w:t></w:r></w:p>'
Replace-ParagraphXml "Also note that the radiation shields" $cutLengthXml

# ---------------------------------------------------------------------
# 13. "*note that the thicknesses were estimated..." paragraph -
#     gramStart/gramEnd split around "bit".
# ---------------------------------------------------------------------
$thicknessXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">*note that the thicknesses were estimated using a micrometer and a drill </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>bit</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> so the accuracy of the measurement isn&#8217;t excellent.</w:t></w:r></w:p>'
Replace-ParagraphXml "note that the thicknesses were estimated" $thicknessXml

# ---------------------------------------------------------------------
# 14. "This calculation is a good starting point..." paragraph -
#     gramStart/gramEnd split around "Thus".
# ---------------------------------------------------------------------
$heatCapacityXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">This calculation is a good starting point but ignores the fact that heat capacity is a function of temperature which trends towards zero at low temperatures, meanwhile cooling rate is also temperature dependent, increasing with temperature. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>Thus</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> this is a gross overestimate of the cooling </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>time</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>.</w:t></w:r></w:p>'
Replace-ParagraphXml "This calculation is a good starting point" $heatCapacityXml

# ---------------------------------------------------------------------
# 15. "Looking at ... Second stage cooling from a Cryomech PT415..." paragraph -
#     spellStart/spellEnd split around "Cryomech".
# ---------------------------------------------------------------------
$cryomechXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">Looking at </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>&#8220;</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">Second stage cooling from a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>Cryomech</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> PT415 cooler at second stage temperatures up to 300 K with cooling on the first-stage from 0 to 250 W</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>&#8221; by green if the second stage is sufficiently cold then the first stage cooling power can be estimated as:</w:t></w:r></w:p>'
Replace-ParagraphXml "Looking at " $cryomechXml
